$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.483333333333333
$ws.Cells.Item(3, 2).Value = 0.45
$ws.Cells.Item(4, 2).Value = 0.533333333333333
$ws.Cells.Item(5, 2).Value = 0.666666666666666
$ws.Cells.Item(6, 2).Value = 0.516666666666666
$ws.Cells.Item(7, 2).Value = 0.483333333333333
$ws.Cells.Item(8, 2).Value = 0.516666666666666
$ws.Cells.Item(9, 2).Value = 0.5
$ws.Cells.Item(10, 2).Value = 0.683333333333333
$ws.Cells.Item(11, 2).Value = 0.5
$ws.Cells.Item(12, 2).Value = 0.45
$ws.Cells.Item(13, 2).Value = 0.666666666666666
$ws.Cells.Item(14, 2).Value = 0.75
$ws.Cells.Item(15, 2).Value = 0.716666666666666
$ws.Cells.Item(16, 2).Value = 0.6
$ws.Cells.Item(17, 2).Value = 0.683333333333333
$ws.Cells.Item(18, 2).Value = 0.816666666666666
$ws.Cells.Item(19, 2).Value = 0.8
$ws.Cells.Item(20, 2).Value = 0.816666666666666
$ws.Cells.Item(21, 2).Value = 0.85
$ws.Cells.Item(22, 2).Value = 0.8
$ws.Cells.Item(23, 2).Value = 0.65
$ws.Cells.Item(24, 2).Value = 0.666666666666666
$ws.Cells.Item(25, 2).Value = 0.916666666666666
$ws.Cells.Item(26, 2).Value = 0.816666666666666
$ws.Cells.Item(27, 2).Value = 0.566666666666666
$ws.Cells.Item(28, 2).Value = 0.783333333333333
$ws.Cells.Item(29, 2).Value = 0.733333333333333
$ws.Cells.Item(30, 2).Value = 0.583333333333333
$ws.Cells.Item(31, 2).Value = 0.766666666666666
$ws.Cells.Item(32, 2).Value = 0.566666666666666
$ws.Cells.Item(33, 2).Value = 0.466666666666666
$ws.Cells.Item(34, 2).Value = 0.433333333333333
$ws.Cells.Item(35, 2).Value = 0.4
$ws.Cells.Item(36, 2).Value = 0.566666666666666
$ws.Cells.Item(37, 2).Value = 0.633333333333333
$ws.Cells.Item(38, 2).Value = 0.65
$ws.Cells.Item(39, 2).Value = 0.866666666666666
$ws.Cells.Item(40, 2).Value = 0.666666666666666
$ws.Cells.Item(41, 2).Value = 0.716666666666666
$ws.Cells.Item(42, 2).Value = 0.649999999999999
$ws.Cells.Item(43, 2).Value = 0.683333333333333
$ws.Cells.Item(44, 2).Value = 0.683333333333333
$ws.Cells.Item(45, 2).Value = 0.783333333333333
$ws.Cells.Item(46, 2).Value = 0.766666666666666
$ws.Cells.Item(47, 2).Value = 0.916666666666666
$ws.Cells.Item(48, 2).Value = 0.783333333333333
$ws.Cells.Item(49, 2).Value = 0.766666666666666
$ws.Cells.Item(50, 2).Value = 0.749999999999999
$ws.Cells.Item(51, 2).Value = 0.783333333333333

$ws.Range("C2:C51").Select()
